# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" (E2) and
# "Correspond Handback DateTime" (H2) timestamps on the zh-cn and de-de
# worksheets to reflect a freshly generated handback status report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-24 13:29:41"
$wsZhCn.Range("H2").Value = "2016-03-24 13:30:09"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-24 13:29:45"
$wsDeDe.Range("H2").Value = "2016-03-24 13:30:30"
